$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExistingRunConfig_Data")

# Final desired values (Level label, B, C, D, E, F) for rows 16-38.
# Rows 1-15 are unchanged. Rows 16-32 are being updated (values shifted/recomputed)
# and rows 33-38 are newly appended (test-case mapping rows for reporting).

$rowData = @{
    16 = @("Level 16", 1, 1, 1, 1, "NO")
    17 = @("Level 17", 1, 1, 1, 1, "NO")
    18 = @("Level 18", 1, 1, 1, 1, "NO")
    19 = @("Level 19", 1, 1, 1, 1, "NO")
    20 = @("Level 20", 1, 1, 1, 1, "NO")
    21 = @("Level 21", 1, 1, 1, 1, "NO")
    22 = @("Level 22", 3, 4, 2, 2, "NO")
    23 = @("Level 23", 1, 24, 1, 1, "NO")
    24 = @("Level 24", 2, 2, 2, 2, "NO")
    25 = @("Level 25", 1, 7, 5, 1, "NO")
    26 = @("Level 26", 1, 3, 2, 4, "NO")
    27 = @("Level 27", 1, 2, 1, 1, "NO")
    28 = @("Level 28", 99, 19, 19, 19, "NO")
    29 = @("Level 28", 99, 19, 19, 7, "NO")
    30 = @("Level 29", 2, 5, 2, 2, "NO")
    31 = @("Level 30", 3, 2, 3, 2, "NO")
    32 = @("Level 31", 11, 11, 11, 11, "NO")
    33 = @("Level 32", 11, 11, 11, 11, "NO")
    34 = @("Level 33", 1, 1, 1, 1, "NO")
    35 = @("Level 34", 3, 7, 5, 7, "NO")
    36 = @("Level 35", 1, 8, 5, 1, "NO")
    37 = @("Level 36", 1, 1, 1, 1, "NO")
    38 = @("Level 37", 2, 1, 1, 1, "NO")
}

foreach ($r in ($rowData.Keys | Sort-Object)) {
    $vals = $rowData[$r]
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}

